$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "IB User with blank Card via EBS (Tagged to Credit Card Brn) at at -12-Dec-2024"
$ws.Range("A3").Value = "Active IB User of at -12-Dec-2024"
